$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header row: "<field>_old" -> "<field>_FV2310" (columns A..J)
#    and "<field>_new" -> "<field>_FV2404" (columns L..U). Column K ("diff")
#    is left untouched.
# ---------------------------------------------------------------------------
$fv2310 = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)
$fv2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $fv2310.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310[$i]
}
for ($i = 0; $i -lt $fv2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2404[$i]
}

# ---------------------------------------------------------------------------
# 2. Freeze the header row (split after row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3. Turn the data range into a native Excel Table ("Table1") so headers get
#    filter buttons and the range is recognised as a ListObject.
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("A1:U84")
$table = $ws.ListObjects.Add(1, $dataRange, [Type]::Missing, 1)
$table.Name = "Table1"
$table.TableStyle = ""
